$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-07 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-08 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("123÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "654÷7=", 2) | Out-Null
$d.Content.Find.Execute("995÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "832÷9=", 2) | Out-Null
$d.Content.Find.Execute("320÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "339÷6=", 2) | Out-Null
$d.Content.Find.Execute("131÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "131÷3=", 2) | Out-Null
$d.Content.Find.Execute("570÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "333÷4=", 2) | Out-Null
$d.Content.Find.Execute("356÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "486÷6=", 2) | Out-Null
$d.Content.Find.Execute("385÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "343÷2=", 2) | Out-Null
$d.Content.Find.Execute("758÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "246÷4=", 2) | Out-Null
$d.Content.Find.Execute("623÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "696÷2=", 2) | Out-Null
$d.Content.Find.Execute("743÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "210÷3=", 2) | Out-Null
$d.Content.Find.Execute("555÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "332÷7=", 2) | Out-Null
$d.Content.Find.Execute("513÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "956÷4=", 2) | Out-Null
$d.Content.Find.Execute("305÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "783÷7=", 2) | Out-Null
$d.Content.Find.Execute("495÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "410÷4=", 2) | Out-Null
$d.Content.Find.Execute("367÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "869÷8=", 2) | Out-Null
$d.Content.Find.Execute("477÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "777÷6=", 2) | Out-Null
$d.Content.Find.Execute("578÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "302÷8=", 2) | Out-Null
$d.Content.Find.Execute("928÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "645÷2=", 2) | Out-Null
$d.Content.Find.Execute("551÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "938÷3=", 2) | Out-Null
$d.Content.Find.Execute("382÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "411÷6=", 2) | Out-Null
$d.Content.Find.Execute("112÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "609÷8=", 2) | Out-Null
$d.Content.Find.Execute("834÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "487÷6=", 2) | Out-Null
$d.Content.Find.Execute("608÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "273÷2=", 2) | Out-Null
$d.Content.Find.Execute("604÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "499÷2=", 2) | Out-Null
$d.Content.Find.Execute("689÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "519÷2=", 2) | Out-Null
